$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (task #3): mark completed (green) and record completion timestamp ---
$ws.Range("A4").Interior.Color = 5296274
$ws.Range("F4").Value = 42274.916666666664

# --- Row 8 (task #7): mark completed (green) and record completion timestamp ---
$ws.Range("A8").Interior.Color = 5296274
$ws.Range("F8").Value = 42274.083333333336

# --- Row 9 (task #8): record completion timestamp ---
$ws.Range("F9").Value = 42274.083333333336

# --- Row 10 (task #9): mark completed (green) and record completion timestamp ---
$ws.Range("A10").Interior.Color = 5296274
$ws.Range("F10").Value = 42274.083333333336

# --- Row 14 (task #13, final assignment): mark completed (orange) and note ---
$ws.Range("A14").Interior.Color = 49407
$ws.Range("F14").Value = 42274.957638888889
$ws.Range("G14").Value = "Не были присланы вовремя все файлы от второй группы"

# --- Update current selection to reflect where the author left off ---
$ws.Range("A14").Select() | Out-Null
